$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H48").Value = 1625
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 1625
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 4875
$ws.Range("N48").Value = -5459

$ws.Range("H56").Value = 1625
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 1625
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 4875
$ws.Range("N56").Value = -5943

$ws.Range("H76").Value = 58539.055
$ws.Range("I76").Value = 58539.055
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 58539.055
$ws.Range("L76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -58224.055

$ws.Range("H79").Value = 58539.055
$ws.Range("I79").Value = 58539.055
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 58539.055
$ws.Range("L79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -57447.055

$ws.Range("H137").Value = 1191.16
$ws.Range("I137").Value = 1093.95
$ws.Range("J137").Value = 1580
$ws.Range("K137").Value = 3281.85
$ws.Range("L137").Value = 4740
$ws.Range("M137").Value = -731.8500000000004
$ws.Range("N137").Value = -9840

$ws.Range("H138").Value = 1399.3096
$ws.Range("I138").Value = 1129.6786
$ws.Range("J138").Value = 1938.5714
$ws.Range("K138").Value = 3389.0358
$ws.Range("L138").Value = 5815.7142
$ws.Range("M138").Value = 1750.9642
$ws.Range("N138").Value = -16095.7142

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3028.5715
$ws.Range("I61").Value = 2675
$ws.Range("J61").Value = 3500
$ws.Range("K61").Value = 2675
$ws.Range("L61").Value = 3500
$ws.Range("M61").Value = -2463
$ws.Range("N61").Value = -3924

$ws.Range("H74").Value = 434.14285
$ws.Range("I74").Value = 389.53845
$ws.Range("J74").Value = 1014
$ws.Range("K74").Value = 389.53845
$ws.Range("L74").Value = 1014
$ws.Range("M74").Value = 484.46155
$ws.Range("N74").Value = -2762

$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("L76").ClearContents()
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = 0

$ws.Range("H77").Value = 434.14285
$ws.Range("I77").Value = 389.53845
$ws.Range("J77").Value = 1014
$ws.Range("K77").Value = 1947.69225
$ws.Range("L77").Value = 5070
$ws.Range("M77").Value = 2420.30775
$ws.Range("N77").Value = -13806

$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("L79").ClearContents()
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = 0

$ws.Range("H132").Value = 3875.0364
$ws.Range("I132").Value = 4240.091
$ws.Range("J132").Value = 3327.4546
$ws.Range("K132").Value = 12720.273
$ws.Range("L132").Value = 9982.363799999999
$ws.Range("M132").Value = -10190.273
$ws.Range("N132").Value = -15042.3638

$ws.Range("H136").Value = 3028.5715
$ws.Range("I136").Value = 2675
$ws.Range("J136").Value = 3500
$ws.Range("K136").Value = 8025
$ws.Range("L136").Value = 10500
$ws.Range("M136").Value = -5475
$ws.Range("N136").Value = -15600

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 308.22223
$ws.Range("I22").Value = 293.2
$ws.Range("J22").Value = 383.33334
$ws.Range("K22").Value = 293.2
$ws.Range("L22").Value = 383.33334
$ws.Range("M22").Value = 56.80000000000001
$ws.Range("N22").Value = -1083.33334

$ws.Range("H99").Value = 2447.55
$ws.Range("I99").Value = 2051.375
$ws.Range("J99").Value = 2711.6667
$ws.Range("K99").Value = 2051.375
$ws.Range("L99").Value = 2711.6667
$ws.Range("M99").Value = -553.375
$ws.Range("N99").Value = -5707.6667

$ws.Range("H122").Value = 1490.5
$ws.Range("I122").Value = 784.8
$ws.Range("J122").Value = 2666.6667
$ws.Range("K122").Value = 2354.4
$ws.Range("L122").Value = 8000.000100000001
$ws.Range("M122").Value = 95.60000000000036
$ws.Range("N122").Value = -12900.0001

$ws.Range("H126").Value = 2447.55
$ws.Range("I126").Value = 2051.375
$ws.Range("J126").Value = 2711.6667
$ws.Range("K126").Value = 6154.125
$ws.Range("L126").Value = 8135.000100000001
$ws.Range("M126").Value = -3684.125
$ws.Range("N126").Value = -13075.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1304.8052
$ws.Range("I68").Value = 1151.6888
$ws.Range("J68").Value = 1520.125
$ws.Range("K68").Value = 3455.0664
$ws.Range("L68").Value = 4560.375
$ws.Range("M68").Value = -2644.0664
$ws.Range("N68").Value = -6182.375

$ws.Range("H71").Value = 1304.8052
$ws.Range("I71").Value = 1151.6888
$ws.Range("J71").Value = 1520.125
$ws.Range("K71").Value = 10365.1992
$ws.Range("L71").Value = 13681.125
$ws.Range("M71").Value = -6309.199199999999
$ws.Range("N71").Value = -21793.125

$ws.Range("H94").Value = 6683.0835
$ws.Range("I94").Value = 2425
$ws.Range("J94").Value = 8812.125
$ws.Range("K94").Value = 7275
$ws.Range("L94").Value = 26436.375
$ws.Range("M94").Value = -6599
$ws.Range("N94").Value = -27788.375

$ws.Range("H97").Value = 339.9091
$ws.Range("I97").Value = 291.2857
$ws.Range("J97").Value = 425
$ws.Range("K97").Value = 873.8571000000001
$ws.Range("L97").Value = 1275
$ws.Range("M97").Value = -377.8571000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 59727.2
$ws.Range("I132").Value = 68898.07000000001
$ws.Range("J132").Value = 4702
$ws.Range("K132").Value = 206694.21
$ws.Range("L132").Value = 14106
$ws.Range("M132").Value = -204164.21
$ws.Range("N132").Value = -19166

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 549.46155
$ws.Range("I22").Value = 510.5
$ws.Range("J22").Value = 566.7778
$ws.Range("K22").Value = 510.5
$ws.Range("L22").Value = 566.7778
$ws.Range("M22").Value = -215.5
$ws.Range("N22").Value = -1156.7778

$ws.Range("H27").Value = 549.46155
$ws.Range("I27").Value = 510.5
$ws.Range("J27").Value = 566.7778
$ws.Range("K27").Value = 510.5
$ws.Range("L27").Value = 566.7778
$ws.Range("M27").Value = -403.5
$ws.Range("N27").Value = -780.7778

$ws.Range("H46").Value = 1655.5
$ws.Range("I46").Value = 890.9091
$ws.Range("J46").Value = 2150.2354
$ws.Range("K46").Value = 890.9091
$ws.Range("L46").Value = 2150.2354
$ws.Range("M46").Value = -702.9091
$ws.Range("N46").Value = -2526.2354

$ws.Range("H68").Value = 1230.4286
$ws.Range("I68").Value = 1230.4286
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 1230.4286
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -481.4286

$ws.Range("H71").Value = 1230.4286
$ws.Range("I71").Value = 1230.4286
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 6152.143
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -2408.143

$ws.Range("H93").Value = 1127371
$ws.Range("I93").Value = 2458211
$ws.Range("J93").Value = 1275.6154
$ws.Range("K93").Value = 2458211
$ws.Range("L93").Value = 1275.6154
$ws.Range("M93").Value = -2456963

$ws.Range("H100").Value = 2519.7
$ws.Range("I100").Value = 2244.7778
$ws.Range("J100").Value = 4994
$ws.Range("K100").Value = 2244.7778
$ws.Range("L100").Value = 4994
$ws.Range("M100").Value = -1703.7778

$ws.Range("H122").Value = 2247.2
$ws.Range("I122").Value = 2628.5
$ws.Range("J122").Value = 1993
$ws.Range("K122").Value = 7885.5
$ws.Range("L122").Value = 5979
$ws.Range("M122").Value = -5435.5
$ws.Range("N122").Value = -10879

$ws.Range("H132").Value = 11488.454
$ws.Range("I132").Value = 12652.944
$ws.Range("J132").Value = 6248.25
$ws.Range("K132").Value = 37958.83199999999
$ws.Range("L132").Value = 18744.75
$ws.Range("M132").Value = -35428.83199999999
$ws.Range("N132").Value = -23804.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4500
$ws.Range("I62").Value = 4000
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 4000
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -3376
$ws.Range("N62").Value = -6248

$ws.Range("H65").Value = 4500
$ws.Range("I65").Value = 4000
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 20000
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -16880
$ws.Range("N65").Value = -31240

$ws.Range("H122").Value = 2103.0557
$ws.Range("I122").Value = 2143.5715
$ws.Range("J122").Value = 1961.25
$ws.Range("K122").Value = 6430.7145
$ws.Range("L122").Value = 5883.75
$ws.Range("M122").Value = -3980.7145
$ws.Range("N122").Value = -10783.75
